$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "acceleration (0-60 mph)" header to "acceleration (0-100)"
$ws.Range("G1").Value = "acceleration (0-100)"

# Reflect the updated selection (G1) as recorded in the saved view state
$ws.Range("G1").Select()
